$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Clear the value in B1 while keeping its existing style/number format
$ws.Range("B1").ClearContents()

# Update the selected cell to C2 (matches the sheetView selection change in the diff)
$ws.Range("C2").Select()
